# Auto-generated edit script: refreshes the "cryptos" price/volume table
# (Price column D, Volume(1h) column E, and one coin swap in row 51)
# to match the values captured in the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.065.80"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "3.419.60"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.74"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.58"
$ws.Range("E6").Value = "  +5.21%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.686"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.48"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.43"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.91"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "3.420.47"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "62.084.81"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.10"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000131"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.47"
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "313.58"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.88"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.16"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.76"
$ws.Range("E25").Value = "  +8.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.65"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.18"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +5.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "42.62"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.34"
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0483"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.37"
$ws.Range("E36").Value = "  -2.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.311"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.11"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.76"
$ws.Range("E45").Value = "  -3.87%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.31"
$ws.Range("E47").Value = "  -4.42%  "
$ws.Range("D48").Value = "2.118.94"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.91"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("B51").Value = "Fetch.AI"
$ws.Range("C51").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.66"
$ws.Range("E51").Value = "  +19.61%  "
